$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data for the "MatplotLib page" (Doubly Linked List sorting) column B
$ws.Range("B2").Value = 0.00099992752075195291
$ws.Range("B3").Value = 0.0019993782043457001

# Update the active selection to C4
$ws.Range("C4").Select()
